# Update "want to go" counts (column F) on the exhibition ("展览"),
# performance ("演出"), and combined ("全部类型") sheets to match the
# freshly scraped totals.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 284
$ws1.Range("F4").Value = 2743
$ws1.Range("F5").Value = 59
$ws1.Range("F6").Value = 579

# 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# 全部类型 (All types - combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2
$ws4.Range("F5").Value = 284
$ws4.Range("F6").Value = 2743
$ws4.Range("F7").Value = 59
$ws4.Range("F8").Value = 579
